$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.641.26'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '1.862.41'
$ws.Range("E3").Value = '  +2.05%  '
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'245.62"
$ws.Range("E5").Value = '  +2.70%  '
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.07741"
$ws.Range("E8").Value = '  +1.86%  '
$ws.Range("D10").Value = "'23.70"
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").Value = "'0.07772"
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'5.161"
$ws.Range("E12").Value = '  +2.60%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.855.10'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = "'92.33"
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("D15").Value = "'0.6928"
$ws.Range("E15").Value = '  +3.43%  '
$ws.Range("E16").Value = '  +3.24%  '
$ws.Range("D17").Value = '29.622.18'
$ws.Range("D18").Value = "'0.000008372"
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '2.108.31'
$ws.Range("E19").Value = '  +1.68%  '
$ws.Range("D20").Value = "'242.04"
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = "'7.620"
$ws.Range("E23").Value = '  +3.56%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = "'0.1510"
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("D26").Value = "'8.920"
$ws.Range("E26").Value = '  +2.60%  '
$ws.Range("D27").Value = "'159.64"
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("D29").Value = "'1.535"
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").Value = "'4.261"
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("D31").Value = "'4.198"
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").Value = "'1.195"
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("D33").Value = "'0.05114"
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("E34").Value = '  +5.67%  '
$ws.Range("E35").Value = '  +5.53%  '
$ws.Range("D36").Value = "'1.159"
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("D38").Value = '1.332.26'
$ws.Range("E38").Value = '  +11.64%  '
$ws.Range("E39").Value = '  +2.83%  '
$ws.Range("D40").Value = "'2.738"
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = "'0.9639"
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("D42").Value = "'5.943"
$ws.Range("E42").Value = '  +14.38%  '
$ws.Range("D43").Value = "'106.53"
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'0.00000000127"
$ws.Range("E45").Value = '  +4.29%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'9.775"
$ws.Range("E46").Value = '  +3.90%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '2.008.30'
$ws.Range("E47").Value = '  +1.56%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = "'0.5214"
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").Value = "'1.791"
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'64.70"
$ws.Range("E50").Value = '  +4.31%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = "'7.007"
$ws.Range("E51").Value = '  +2.40%  '